$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy an existing data row (row 2) to rows 6 and 7 so formatting/styles
# (e.g. the date number format on column B) are preserved exactly.
$ws.Range("A2:D2").Copy($ws.Range("A6:D6"))
$ws.Range("A2:D2").Copy($ws.Range("A7:D7"))

# New row 6: d
$ws.Range("A6").Value = "d"

# New row 7: e
$ws.Range("A7").Value = "e"

# Update selection to F10 as shown in the diff
$ws.Range("F10").Select()
